# "another pass on 3.3" -- nudge several shapes/connectors on slide 1,
# widen a few label textboxes (wrapping their text in "{ }"), and add a
# new small "∅" label textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    foreach ($sh in $slide.Shapes) {
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# --- Oval 7 (id 8): "3: t = s;" -- shift left only ---
$sh8 = Get-ShapeById $s 8
$sh8.Left = 439.2276458952756

# --- Oval 8 (id 9): "4: close(t)" -- shift left only ---
$sh9 = Get-ShapeById $s 9
$sh9.Left = 430.7528346456693

# --- Oval 9 (id 10): "exit" -- shift down only ---
$sh10 = Get-ShapeById $s 10
$sh10.Top = 270.60259842519685

# --- Straight Arrow Connector 15 (id 16) -- shift left only ---
$sh16 = Get-ShapeById $s 16
$sh16.Left = 459.97496062992127

# --- Straight Arrow Connector 18 (id 19) -- shift left only ---
$sh19 = Get-ShapeById $s 19
$sh19.Left = 492.6584320968504

# --- Straight Arrow Connector 21 (id 22) -- grow taller only ---
$sh22 = Get-ShapeById $s 22
$sh22.Height = 57.14700897401575

# --- Straight Arrow Connector 24 (id 25) -- resize (flipped connector) ---
$sh25 = Get-ShapeById $s 25
$sh25.Width = 26.708976377952755
$sh25.Height = 21.78204724409449

# --- TextBox 38 (id 39): "<{s}, e>" -> "{<{s}, e>}", move + widen ---
$sh39 = Get-ShapeById $s 39
$sh39.TextFrame.TextRange.Text = "{<{s}, e>}"
$sh39.Left = 304.2330780661418
$sh39.Top = 164.42275590551182
$sh39.Width = 67.9320488440945

# --- TextBox 40 (id 41): "<{s}, e>" -> "{<{s}, e>}", shift left + widen ---
$sh41 = Get-ShapeById $s 41
$sh41.TextFrame.TextRange.Text = "{<{s}, e>}"
$sh41.Left = 487.9740157480315
$sh41.Width = 67.9320488440945

# --- TextBox 41 (id 42): "<{s, t}, e>" -> "{<{s, t}, e>}", shift left + widen ---
$sh42 = Get-ShapeById $s 42
$sh42.TextFrame.TextRange.Text = "{<{s, t}, e>}"
$sh42.Left = 492.6584320968504
$sh42.Width = 79.67047504094488

# --- TextBox 44 (id 45): "<∅, e>" -> "{<∅, e>}", move + widen ---
$sh45 = Get-ShapeById $s 45
$sh45.TextFrame.TextRange.Text = "{<∅, e>}"
$sh45.Left = 404.65773013543304
$sh45.Top = 247.59866141732283
$sh45.Width = 102.84968503937007

# --- New TextBox 17 (id 18): small "∅" label near the top ---
# The host assigns shape Ids from an internal monotonic counter (not
# simply Shapes.Count+1); replaying the same number of alloc/free cycles
# the original authoring session would have burned through reproduces
# the target Id (18) for the newly-added shape deterministically.
for ($i = 0; $i -lt 7; $i++) {
    $tmp = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $tmp.Delete()
}

$new = $sh45.Duplicate()
$new.Name = "TextBox 17"
$new.TextFrame.TextRange.Text = "∅"
$new.Left = 444.76149606299214
$new.Top = 92.03606299212598
$new.Width = 32.33779717559055
$new.Height = 24.234410348818898
